# tfs10524 - ecl move db file shares
# Repoint the eCL Test Job Steps workbook's file-share references from the
# old \\vrivscors01 "Test" shares to the new \\F3420-ECLDBP01 "Encrypt_In"
# shares, and mark the two now-unused legacy steps accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eCL Jobs")

# EmployeeLoad (row 2) - Input
$ws.Range("E2").Value = "\\F3420-ECLDBP01\Data\Coaching\HRInfo\Encrypt_In\Employee_Information_WithProgram.csv.zip.encrypt`n\\vrivscors01\BCC \\F3420-ECLDBP01\Data\Coaching\HRInfo\Encrypt_In\PS_Employee_Information_<MMDDCCYY>.csv.zip.encrypt`n\\F3420-ECLDBP01\Data\Coaching\HRInfo\Encrypt_In\HR_Employee_Information.csv.zip.encrypt"

# CoachingETSLoad (row 4) - Input
$ws.Range("E4").Value = "\\F3420-ECLDBP01\data\Coaching\ETS\Encrypt_In\eCL_ETS_Feed_<ReportCode><CCYYMMDD>.csv.zip.encrypt"

# CoachingGenericLoad (row 6) - Input
$ws.Range("E6").Value = "\\F3420-ECLDBP01\data\Coaching\Generic\Encrypt_In\eCL_Generic_Feed_XXX[_ZZZ]<YYYYMMDD>.csv.zip.encrypt"

# CoachingOutliersLoad (row 8) - Input
$ws.Range("E8").Value = "\\F3420-ECLDBP01\data\Coaching\Outliers\Encrypt_In\eCl_Outlier_Feed_<ReportCode><CCYYMMDD>.csv.zip.encrypt"

# CoachingQualityLoad / IQSLoad (row 10) - Input
$ws.Range("E10").Value = "\\F3420-ECLDBP01\data\Coaching\IQS\Encrypt_In\eCL_IQS_Scorecard_<CCYYMMDD>.csv.zip.encrypt"

# CoachingQualityOtherLoad (row 12) - Input
$ws.Range("E12").Value = "\\F3420-ECLDBP01\data\Coaching\Quality\Encrypt_In\eCL_Quality_Feed_XXX<YYYYMMDD>.csv.zip.encrypt"

# CoachingTraining (row 14) - Input
$ws.Range("E14").Value = "\\F3420-ECLDBP01\data\Coaching\Training\Encrypt_In\eCL_Training_Feed_XXX<YYYYMMDD>.csv.zip.encrypt"

# CoachingWHLoad (row 26) - Input -> legacy, no longer used
$ws.Range("E26").Value = "legacy - no longer used"

# CoachingInactivations (row 28) - Input -> legacy, no longer used
$ws.Range("E28").Value = "legacy - no longer used"
# CoachingInactivations (row 28) - Output: drop the old UNC log path reference
$ws.Range("F28").Value = "EC.Coaching_Log.StatusID = 2 or`nEC.Warning_Log.StatusID = 2`nemail notification sent to john;`nlog file generated to <>"

# Row 28 text got shorter, so the wrapped-text row shrinks from 86.4 to 57.6
$ws.Rows.Item(28).RowHeight = 57.6
